$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new headers (I0, IF) in columns I and J, row 1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the rest of row 1 (copy style from H1).
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF), rows 2-24.
$data = @{
    2  = @(7, 7)
    3  = @(8, 9)
    4  = @(9, 9)
    5  = @(6, 7)
    6  = @(9, 9)
    7  = @(9, 9)
    8  = @(7, 8)
    9  = @(4, 6)
    10 = @(8, 8)
    11 = @(11, 11)
    12 = @(6, 6)
    13 = @(7, 8)
    14 = @(7, 7)
    15 = @(6, 6)
    16 = @(7, 7)
    17 = @(6, 7)
    18 = @(5, 7)
    19 = @(6, 8)
    20 = @(7, 8)
    21 = @(5, 5)
    22 = @(1, 2)
    23 = @(5, 6)
    24 = @(3, 4)
}

foreach ($r in $data.Keys) {
    $pair = $data[$r]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}
